# work-pulse.xlsx update
# - rename "Лист1" -> "Пульс"
# - add new sheet "Начисления" with a small lookup table
# - add today's (row 10 / 2024-06-10) tracking numbers and fix row 9 "отдых" value
# - re-point the chart series formulas at the renamed sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- rename the main sheet ---
$ws.Name = "Пульс"

# --- data fixes on the existing rows ---
$ws.Range("E9").Value = 2

# --- fill in the previously-empty row 10 (2024-06-10) ---
$ws.Range("B10").Value = 3
$ws.Range("C10").Value = 3
$ws.Range("D10").Value = 3
$ws.Range("E10").Value = 1

# --- point chart series formulas at the renamed sheet and refresh their cache ---
$chartObj = $ws.ChartObjects(1)
$chart = $chartObj.Chart
$cols = @("B", "C", "D", "E")
for ($i = 1; $i -le 4; $i++) {
    $col = $cols[$i - 1]
    $ser = $chart.SeriesCollection($i)
    $ser.Formula = "=SERIES(Пульс!`$$col`$1,Пульс!`$A`$2:`$A`$19,Пульс!`$$col`$2:`$$col`$19,$i)"
}
$chart.Refresh()

# --- add the new "Начисления" sheet (lookup of habit -> category -> points) ---
$ws2 = $wb.Worksheets.Add($null, $ws)
$ws2.Name = "Начисления"

$ws2.Range("A2").Value = 'привычка "закинуть чтото в блог"'
$ws2.Range("B2").Value = "потик.ио"
$ws2.Range("C2").Value = 1

$ws2.Range("A3").Value = "чтение"
$ws2.Range("B3").Value = "отдых"
$ws2.Range("C3").Value = 1

$ws2.Range("A1").Value = "Что"
$ws2.Range("B1").Value = "Категория"
$ws2.Range("C1").Value = "Сколько баллов"

$ws2.Columns.AutoFit()
$ws2.Range("C9").Select()

# --- re-activate the main sheet and restore the expected selection ---
$ws.Activate()
$ws.Range("F21").Select()
